$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.087.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "'1.891.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'314.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.5041"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.82%  "
$ws.Range("D8").Value = "'0.3899"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").Value = "'0.09241"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.74%  "
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").Value = "'41.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "'6.382"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("D13").Value = "'20.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "'1.890.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "'7.295"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'92.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "'0.06645"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "'17.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'6.205"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").Value = "'28.138.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "'11.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'2.318"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "'2.105.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").Value = "'2.540"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("D28").Value = "'158.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "'20.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").Value = "'126.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").Value = "'5.601"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "'3.600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").Value = "'0.06614"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "'1.341"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.72%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "'0.2196"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").Value = "'1.218"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6445"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").Value = "'4.969"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'13.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "'0.6051"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "'1.302"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Value = "'3.689"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "'2.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").Value = "'121.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("D51").Value = "'1.195"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.55%  "
